$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '28.014.77'
$ws.Cells.Item(2, 5).Value = '  +1.80%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.651.89'
$ws.Cells.Item(3, 5).Value = '  +2.10%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.17%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '213.78'
$ws.Cells.Item(5, 5).Value = '  +1.26%  '

$ws.Cells.Item(6, 5).Value = '  +0.83%  '

$ws.Cells.Item(7, 5).Value = '  -0.17%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '23.60'
$ws.Cells.Item(8, 5).Value = '  +3.72%  '

$ws.Cells.Item(9, 5).Value = '  +1.63%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.885.30'
$ws.Cells.Item(12, 5).Value = '  +2.03%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '1.644.16'
$ws.Cells.Item(13, 5).Value = '  +1.65%  '

$ws.Cells.Item(14, 5).Value = '  +1.67%  '

$ws.Cells.Item(15, 5).Value = '  +2.65%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '65.74'
$ws.Cells.Item(16, 5).Value = '  +1.16%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '28.000.22'
$ws.Cells.Item(17, 5).Value = '  +1.71%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '233.01'
$ws.Cells.Item(18, 5).Value = '  +1.40%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '7.69'
$ws.Cells.Item(19, 5).Value = '  +2.11%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.0₃0723'
$ws.Cells.Item(20, 5).Value = '  +0.50%  '

$ws.Cells.Item(21, 5).Value = '  -0.15%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '10.69'
$ws.Cells.Item(22, 5).Value = '  +5.33%  '

$ws.Cells.Item(23, 5).Value = '  +2.81%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.16'
$ws.Cells.Item(24, 5).Value = '  +3.56%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '152.35'
$ws.Cells.Item(25, 5).Value = '  +1.88%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '6.91'
$ws.Cells.Item(26, 5).Value = '  +1.42%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '15.75'
$ws.Cells.Item(27, 5).Value = '  +1.22%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.112'
$ws.Cells.Item(28, 5).Value = '  +0.65%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.20'
$ws.Cells.Item(30, 5).Value = '  +1.79%  '

$ws.Cells.Item(31, 5).Value = '  +0.55%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.35'
$ws.Cells.Item(32, 5).Value = '  +2.68%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.454.36'
$ws.Cells.Item(33, 5).Value = '  +1.00%  '

$ws.Cells.Item(34, 5).Value = '  +1.50%  '

$ws.Cells.Item(35, 5).Value = '  +2.46%  '

$ws.Cells.Item(36, 5).Value = '  -0.57%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.892'
$ws.Cells.Item(37, 5).Value = '  +3.59%  '

$ws.Cells.Item(38, 5).Value = '  +0.86%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.562'
$ws.Cells.Item(39, 5).Value = '  +0.25%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.919'
$ws.Cells.Item(40, 5).Value = '  -1.96%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '69.46'
$ws.Cells.Item(41, 5).Value = '  +2.42%  '

$ws.Cells.Item(42, 5).Value = '  +2.68%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.999'
$ws.Cells.Item(43, 5).Value = '  -0.17%  '

$ws.Cells.Item(44, 5).Value = '  +0.21%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2.24'
$ws.Cells.Item(45, 5).Value = '  +1.27%  '

$ws.Cells.Item(46, 5).Value = '  +6.37%  '

$ws.Cells.Item(47, 5).Value = '  -0.86%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.794.02'
$ws.Cells.Item(48, 5).Value = '  +1.83%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '88.79'
$ws.Cells.Item(49, 5).Value = '  +2.83%  '

$ws.Cells.Item(50, 5).Value = '  +0.19%  '

$ws.Cells.Item(51, 5).Value = '  +1.10%  '
